# Updated CVDs for the month
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Milwaukee Pmc Hq Wisconsin
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Milwaukee Pmc Hq Wisconsin")
$ws.Range("E2").Value = 0.1163
$ws.Range("E3").Value = 0.1163
$ws.Range("E4").Value = 0.1163
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 0
$ws.Range("V4").Value = 0
$ws.Range("W4").Value = 0
$ws.Range("O7").Value = ""

# ---------------------------------------------------------------------
# 2. Monterrey Rbm Mexico
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Monterrey Rbm Mexico")
$ws.Range("O4").Value = ""
$ws.Range("P4").Value = 0.0833333333333333
$ws.Range("Q4").Value = 0.0833333333333333
$ws.Range("R4").Value = 0.25
$ws.Range("S4").Value = 0.0833333333333333
$ws.Range("T4").Value = 0.0833333333333333
$ws.Range("U4").Value = 0.0833333333333333
$ws.Range("V4").Value = 0.25
$ws.Range("W4").Value = 1

# ---------------------------------------------------------------------
# 3. Piedras Negras Fasco Mexico
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Piedras Negras Fasco Mexico")
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0.0833333333333333
$ws.Range("Q4").Value = 0.0833333333333333
$ws.Range("R4").Value = 0.25
$ws.Range("S4").Value = 0.0833333333333333
$ws.Range("T4").Value = 0.0833333333333333
$ws.Range("U4").Value = 0.0833333333333333
$ws.Range("V4").Value = 0.25
$ws.Range("W4").Value = 1

# ---------------------------------------------------------------------
# 4. Rosemont Illinois
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Rosemont Illinois")
$ws.Range("O7").Value = ""

# ---------------------------------------------------------------------
# 5. Apodaca Pmc Plant 1 Mexico
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Apodaca Pmc Plant 1 Mexico")
$ws.Range("O3").Value = ""

# ---------------------------------------------------------------------
# 6. Braintree Massachusetts
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Braintree Massachusetts")
$ws.Range("E2").Value = 0.2564
$ws.Range("E3").Value = 0.2564
$ws.Range("E4").Value = 0.2564
$ws.Range("O4").Value = 0.3333

# ---------------------------------------------------------------------
# 7. El Paso Texas - EPC
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("El Paso Texas - EPC")
$ws.Range("O4").Value = ""

# ---------------------------------------------------------------------
# 8. Faridabad India
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Faridabad India")
$ws.Range("O5").Value = ""

# ---------------------------------------------------------------------
# 9. Fort Wayne Indiana
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Fort Wayne Indiana")
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0.0833333333333333
$ws.Range("Q4").Value = 0.0833333333333333
$ws.Range("R4").Value = 0.25
$ws.Range("S4").Value = 0.0833333333333333
$ws.Range("T4").Value = 0.0833333333333333
$ws.Range("U4").Value = 0.0833333333333333
$ws.Range("V4").Value = 0.25
$ws.Range("W4").Value = 1

# ---------------------------------------------------------------------
# 10. Juarez Casa SS
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Juarez Casa SS")
$ws.Range("E5").Value = 0
$ws.Range("E6").Value = 0

# New row 7 (AOP/Commit-Forecast row for Internal Fill Rate), mirroring
# the layout of rows 2-6.
$ws.Range("A7").Value = "Corporate"
$ws.Range("B7").Value = "Corp Legal"
$ws.Range("C7").Value = "Juarez Casa SS"
$ws.Range("D7").Value = "Internal Fill Rate"
$ws.Range("E7").NumberFormat = "0.0%"
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = "Commit/Forecast"
$ws.Range("G7:N7").NumberFormat = "0.0%"
$ws.Range("O7:W7").NumberFormat = "0.0%"
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = 0
$ws.Range("Q7").Value = 0
$ws.Range("R7").Value = 0
$ws.Range("S7").Value = 0
$ws.Range("T7").Value = 0
$ws.Range("U7").Value = 0
$ws.Range("V7").Value = 0
$ws.Range("W7").Value = 0

# ---------------------------------------------------------------------
# 11. Juarez Mej SS
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Juarez Mej SS")
$ws.Range("E2").Value = 0.7143
$ws.Range("E3").Value = 0.7143
$ws.Range("E4").Value = 0.7143
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0.0833333333333333
$ws.Range("Q4").Value = 0.0833333333333333
$ws.Range("R4").Value = 0.25
$ws.Range("S4").Value = 0.0833333333333333
$ws.Range("T4").Value = 0.0833333333333333
$ws.Range("U4").Value = 0.0833333333333333
$ws.Range("V4").Value = 0.25
$ws.Range("W4").Value = 1

Write-Host "Edits applied"
